$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.372438311576843
$ws.Range("B1").Value = 1.74567437171936
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.874284625053406
$ws.Range("E1").Value = 0.7593633532524109
